$d = $word.ActiveDocument

$replacements = @(
    @{old = "871÷6="; new = "229÷8="},
    @{old = "707÷6="; new = "230÷2="},
    @{old = "874÷2="; new = "157÷7="},
    @{old = "450÷7="; new = "744÷8="},
    @{old = "343÷5="; new = "817÷8="},
    @{old = "401÷7="; new = "771÷9="},
    @{old = "568÷2="; new = "466÷9="},
    @{old = "548÷6="; new = "872÷4="},
    @{old = "145÷2="; new = "373÷3="},
    @{old = "198÷9="; new = "493÷7="},
    @{old = "175÷5="; new = "952÷8="},
    @{old = "317÷6="; new = "783÷2="},
    @{old = "313÷7="; new = "781÷3="},
    @{old = "890÷3="; new = "580÷5="},
    @{old = "900÷3="; new = "980÷2="},
    @{old = "894÷8="; new = "952÷6="},
    @{old = "209÷5="; new = "333÷5="},
    @{old = "518÷9="; new = "635÷4="},
    @{old = "899÷5="; new = "436÷8="},
    @{old = "194÷7="; new = "450÷3="},
    @{old = "429÷8="; new = "641÷9="},
    @{old = "610÷5="; new = "639÷8="},
    @{old = "721÷9="; new = "655÷5="},
    @{old = "577÷2="; new = "545÷3="},
    @{old = "966÷5="; new = "185÷3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
